$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = " (e.g. when the minimum number appears first, last or anywhere in between)"
$find.Replacement.ClearFormatting()
$find.Replacement.Text = ""
$find.Execute(
    " (e.g. when the minimum number appears first, last or anywhere in between)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2
)
